# Update "average_county_temperature" column (AD) values on Sheet1
# using refreshed NOAA temperature data, per the commit message:
#   "Updated temperature with NOAA data"
#
# Each contiguous block of rows below corresponds to a distinct
# facility/location whose average_county_temperature value was
# recalculated from 10 (placeholder) to the NOAA-derived figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AD30:AD33").Value = 19.79629629629628
$ws.Range("AD34:AD37").Value = 16.86342592592595
$ws.Range("AD38:AD41").Value = 5.486111111111112
$ws.Range("AD42:AD57").Value = 14.96875
$ws.Range("AD58:AD61").Value = 17.25771604938272
$ws.Range("AD90:AD105").Value = 13.75752314814816
$ws.Range("AD110:AD121").Value = 14.96875
$ws.Range("AD142:AD157").Value = 13.0158303464755
$ws.Range("AD158:AD161").Value = 16.86342592592595
$ws.Range("AD162:AD169").Value = -3.847222222222223
$ws.Range("AD178:AD181").Value = 12.41429539295394
$ws.Range("AD182:AD185").Value = 19.60879629629628
